# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Updates Price (D) and Volume(1h) (E) text values for the affected rows, and
# re-applies the new coin ranking order for the rows whose rank order changed
# (Coin/Link/Price/Volume all move together in those cases).
#
# Every cell in this sheet is stored as text (inline strings in the source
# workbook), including numeric-looking prices like "0.160" or "0.0000252".
# Setting NumberFormat to "@" (Text) before writing the value keeps Excel from
# re-interpreting these as numbers (which would drop significant trailing/
# leading zeros or switch to scientific notation). Resetting the style back to
# "Normal" afterwards avoids leaving a stray text-format style on the cell so
# the cell keeps its original (unstyled) appearance.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2 (columns D, E)
Set-TextValue "D2" "67.760.05"
Set-TextValue "E2" "  +0.11%  "

# Row 3 (columns D, E)
Set-TextValue "D3" "3.816.20"
Set-TextValue "E3" "  +1.12%  "

# Row 4 (columns E)
Set-TextValue "E4" "  +0.07%  "

# Row 5 (columns D, E)
Set-TextValue "D5" "603.63"
Set-TextValue "E5" "  +1.47%  "

# Row 6 (columns D, E)
Set-TextValue "D6" "166.36"
Set-TextValue "E6" "  -0.40%  "

# Row 7 (columns E)
Set-TextValue "E7" "  +0.06%  "

# Row 8 (columns E)
Set-TextValue "E8" "  +0.04%  "

# Row 9 (columns D, E)
Set-TextValue "D9" "0.160"
Set-TextValue "E9" "  +0.48%  "

# Row 10 (columns E)
Set-TextValue "E10" "  +1.24%  "

# Row 11 (columns D, E)
Set-TextValue "D11" "6.29"
Set-TextValue "E11" "  -0.52%  "

# Row 12 (columns D, E)
Set-TextValue "D12" "0.0000252"
Set-TextValue "E12" "  -0.79%  "

# Row 13 (columns D, E)
Set-TextValue "D13" "35.96"
Set-TextValue "E13" "  -0.17%  "

# Row 14 (columns D, E)
Set-TextValue "D14" "4.461.14"
Set-TextValue "E14" "  +1.33%  "

# Row 15 (columns D, E)
Set-TextValue "D15" "3.828.89"
Set-TextValue "E15" "  +1.96%  "

# Row 16 (columns D, E)
Set-TextValue "D16" "18.48"
Set-TextValue "E16" "  +1.05%  "

# Row 17 (columns D, E)
Set-TextValue "D17" "67.797.22"
Set-TextValue "E17" "  +0.28%  "

# Row 18 (columns E)
Set-TextValue "E18" "  +1.29%  "

# Row 19 (columns E)
Set-TextValue "E19" "  +1.40%  "

# Row 20 (columns D, E)
Set-TextValue "D20" "462.45"
Set-TextValue "E20" "  +1.33%  "

# Row 21 (columns D, E)
Set-TextValue "D21" "9.88"
Set-TextValue "E21" "  -1.34%  "

# Row 22 (columns E)
Set-TextValue "E22" "  +1.14%  "

# Row 23 (columns E)
Set-TextValue "E23" "  -2.77%  "

# Row 24 (columns D, E)
Set-TextValue "D24" "83.27"
Set-TextValue "E24" "  +0.17%  "

# Row 25 (columns E)
Set-TextValue "E25" "  +2.18%  "

# Row 26 (columns D, E)
Set-TextValue "D26" "2.13"
Set-TextValue "E26" "  -0.02%  "

# Row 27 (columns D, E)
Set-TextValue "D27" "10.07"
Set-TextValue "E27" "  +0.05%  "

# Row 28 (columns E)
Set-TextValue "E28" "  -0.22%  "

# Row 29 (columns D, E)
Set-TextValue "D29" "3.968.29"
Set-TextValue "E29" "  +1.31%  "

# Row 30 (columns E)
Set-TextValue "E30" "  +1.08%  "

# Row 31 (columns E)
Set-TextValue "E31" "  +1.65%  "

# Row 32 (columns D, E)
Set-TextValue "D32" "2.23"
Set-TextValue "E32" "  +0.98%  "

# Row 33 (columns D, E)
Set-TextValue "D33" "29.53"
Set-TextValue "E33" "  -0.39%  "

# Row 34 (columns B, C, D, E)
Set-TextValue "B34" "Binance-PegBSC-USD"
Set-TextValue "C34" "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue "D34" "1.00"
Set-TextValue "E34" "  +0.23%  "

# Row 35 (columns B, C, D, E)
Set-TextValue "B35" "Aptos"
Set-TextValue "C35" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D35" "9.11"
Set-TextValue "E35" "  -0.35%  "

# Row 36 (columns D, E)
Set-TextValue "D36" "0.0999"
Set-TextValue "E36" "  -0.40%  "

# Row 37 (columns D, E)
Set-TextValue "D37" "3.29"
Set-TextValue "E37" "  -0.81%  "

# Row 38 (columns E)
Set-TextValue "E38" "  +0.04%  "

# Row 39 (columns B, C, D, E)
Set-TextValue "B39" "Mantle"
Set-TextValue "C39" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D39" "0.999"
Set-TextValue "E39" "  +0.63%  "

# Row 40 (columns B, C, D, E)
Set-TextValue "B40" "Filecoin"
Set-TextValue "C40" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D40" "5.81"
Set-TextValue "E40" "  +1.12%  "

# Row 41 (columns D, E)
Set-TextValue "D41" "0.999"
Set-TextValue "E41" "  +0.02%  "

# Row 43 (columns D, E)
Set-TextValue "D43" "44.24"
Set-TextValue "E43" "  -3.68%  "

# Row 44 (columns B, C, D, E)
Set-TextValue "B44" "OKB"
Set-TextValue "C44" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D44" "47.69"
Set-TextValue "E44" "  -0.93%  "

# Row 45 (columns B, C, D, E)
Set-TextValue "B45" "TheGraph"
Set-TextValue "C45" "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue "D45" "0.301"
Set-TextValue "E45" "  +0.75%  "

# Row 46 (columns E)
Set-TextValue "E46" "  +15.58%  "

# Row 47 (columns B, C, D, E)
Set-TextValue "B47" "EnergySwap"
Set-TextValue "C47" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D47" "28.35"
Set-TextValue "E47" "  +8.64%  "

# Row 48 (columns D, E)
Set-TextValue "D48" "150.88"
Set-TextValue "E48" "  +1.22%  "

# Row 49 (columns E)
Set-TextValue "E49" "  +0.51%  "

# Row 50 (columns E)
Set-TextValue "E50" "  +2.17%  "

# Row 51 (columns D, E)
Set-TextValue "D51" "390.28"
Set-TextValue "E51" "  +0.39%  "
